$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Group <X>" -> "Group <2>" with a "_GoBack" bookmark placed just before
#    the closing ">" (matches the relocation of the bookmark away from the
#    "Test plays" paragraph later in this script).
# ---------------------------------------------------------------------------
$grp = $d.Content
$grp.Find.Execute(" <X>")
if ($grp.Find.Found) {
    $gs = $grp.Start
    # " <X>" -> chars: [0]=' ' [1]='<' [2]='X' [3]='>'
    $xChar = $d.Range($gs + 2, $gs + 3)
    $xChar.Text = "2"
    $bmPoint = $d.Range($gs + 3, $gs + 3)
    $d.Bookmarks.Add("_GoBack", $bmPoint)
}

# ---------------------------------------------------------------------------
# 2) "Version 1.0" -> "Version 3.0"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Version 1.0", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Version 3.0", 2)

# ---------------------------------------------------------------------------
# 3) Expand the functional-testing paragraph and rewrite the "Test plays"
#    paragraph that follows the "Execution-based Non-Functional Testing"
#    heading.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("various things.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ("all objects. Since the game progresses in a linear fashion, " + `
                         "it is unlikely for a user to misuse functions or to " + [char]0x201C + "break" + [char]0x201D + " the game."), 2)

$d.Content.Find.Execute("are done throughout process, which tests how smoothly the game runs. ", `
                         $true, $false, $false, $false, $false, $true, 1, $false, `
                         "allow us to observe how smoothly the game runs. If it underperforms, we will adjust its traits to reduce any lag or desync.", 2)

# ---------------------------------------------------------------------------
# 4) "through the code written thus far as a group" -> drop "thus far "
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("written thus far as a group", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "written as a group", 2)
